# Applies the "Add files via upload" commit:
#   1. Removes the two review comments (slide 3 and slide 9) - this also
#      drops ppt/comments/comment1.xml & comment2.xml (and their slide
#      relationships / content-type overrides) from the saved package.
#   2. Clears the (now resolved) strikethrough formatting that was left on
#      the "Choix: ..." paragraph of slide 3.

$p = $ppt.ActivePresentation

# --- 1. Delete the lingering review comments ---------------------------
for ($i = 1; $i -le $p.Slides.Count; $i++) {
    $slide = $p.Slides.Item($i)
    for ($c = $slide.Comments.Count; $c -ge 1; $c--) {
        $slide.Comments.Item($c).Delete()
    }
}

# --- 2. Remove the strikethrough from the "Choix:" bullet on slide 3 ---
$slide3 = $p.Slides.Item(3)

for ($sIdx = 1; $sIdx -le $slide3.Shapes.Count; $sIdx++) {
    $shape = $slide3.Shapes.Item($sIdx)
    if (-not $shape.HasTextFrame) { continue }

    $textRange = $shape.TextFrame.TextRange
    $paraCount = $textRange.Paragraphs().Count

    for ($pIdx = 1; $pIdx -le $paraCount; $pIdx++) {
        $paragraph = $textRange.Paragraphs($pIdx, 1)
        if ($paragraph.Text -notlike "Choix:*") { continue }

        $runCount = $paragraph.Runs().Count
        for ($r = 1; $r -le $runCount; $r++) {
            $paragraph.Runs($r, 1).Font.Strikethrough = 0
        }
    }
}
